# Mise à jour du fichier via Shiny
# Refresh the "pro", "ind" and "conso" data sheets with newly generated
# values (the "VA" sheet recomputes automatically since B2:B26 there is
# the formula =pro!Bn -conso!Bn).

$wb = $excel.ActiveWorkbook

$proVals = @(2813111.8068729583,3287367.1300997585,3662740.8201511893,3227953.8746803799,3312469.0678267083,2956003.0162836676,2856163.3978679809,2850467.8224011539,2604093.326811797,2535447.5307234181,2798172.6506881271,2611642.9799922295,2538626.8425992075,2634261.8727078312,2829142.0365330572,3051374.972308469,3202960.7410026803,3610021.6490008826,4151259.4542214484,4342123,4198234,4327768.9849505639,4576163.3169878144,4846140.9448465919,4993665.023576268)
$indVals = @(673991.94009069167,718269.73282074661,537534.72274806001,522163.91703414213,787618.51711440575,839360.98965265614,628156.32657470903,610194.19604593632,877554.0969938268,935204.87307117076,699883.44098080357,679870.27357329766,594625.20264260529,633689.01024709002,474236.670246372,460675.8723810612,820401.82962782832,874297.99663945602,654302.29027910368,635592.5159027878,759411.85557014553,809301.29598861211,605660.41958982416,588341.5595345482,720395.90765966917,767722.20161897468,574543.68733336939,558114.61025531346,687663.8361698807,732839.80192666303,548438.59033344733,532755.98851949605,628226.98533439916,669498.25667975214,501035.39567422151,486708.28823950031,611666.4646001904,651849.79518400761,487827.73784936871,473878.30339793669,756038.2212431404,805706.03129395295,602969.81531883555,585727.8930942748,748005.65525200823,797145.76716447505,596563.53228165582,579504.79772432055,662264.69511116785,705772.0683575297,528181.78986607539,513078.43127343233,750358.33004527888,799653.00055208697,598439.8817118092,581327.49294141249,760520.87391476124,810483.17111604044,606544.90474897204,589200.75284004875,845561.24187892547,901110.25245529122,674367.89772105298,655084.34728810389,887474.45769907639,945776.9502866437,707795.31354791403,687555.90613968635,955342.85605200764,1018103.9522169227,761922.99449182686,740135.80601509369,983492.98152937263,1048103.395686197,784373.81176979362,761944.64215975988,1099515.2516971482,1171747.7303401865,876906.07383044087,851830.94413222466,996468.79170146713,1061931.64963266,794722.5238365134,771997.43279882998,1173486.0635956496,1250578.0428981089,935900.66634731786,909138.58624129777,1249206.7156522898,1331273.1510835534,996290.82428317016,967801.88757536758,1255259.7542722879,1337723.8431077853,1001118.3574372702,972491.37741607951,1334414.932188476,1422079.1077799248,1064247.6830038731,1033815.5198808836)
$consoVals = @(699386.95332367369,815838.87216876808,910544.60749126493,802216.67629314505,821233.59100075858,733988.41856837063,708529.99493727344,705777.2005003267,642461.83180849627,622564.84241331136,686706.02782055223,641265.81947924884,623631.22572862776,646897.48713822605,695326.34948867769,747563.11027763365,785418.86748075124,886298.06879296096,1022894.6740994203,1067530,1142948,1274179.5840800432,1262130.9855690033,1336592.2986663848,1377280.250119308)

$wsPro = $wb.Worksheets.Item("pro")
for ($i = 0; $i -lt $proVals.Length; $i++) {
    $wsPro.Cells.Item($i + 2, 2).Value = $proVals[$i]
}

$wsInd = $wb.Worksheets.Item("ind")
for ($i = 0; $i -lt $indVals.Length; $i++) {
    $wsInd.Cells.Item($i + 2, 2).Value = $indVals[$i]
}

$wsConso = $wb.Worksheets.Item("conso")
for ($i = 0; $i -lt $consoVals.Length; $i++) {
    $wsConso.Cells.Item($i + 2, 2).Value = $consoVals[$i]
}

$wb.Worksheets.Item("VA").Calculate()

# Column B on "pro" narrows slightly (content re-fit) after the refresh.
$wsPro.Columns.Item(2).ColumnWidth = 13.8

# Selection moves to B2:B26 (whole data column) on every refreshed sheet,
# and the "ind" view no longer needs to be scrolled to the bottom.
[void]$wsInd.Range("B2:B26").Select()
[void]$wb.Worksheets.Item("VA").Range("B2:B26").Select()
[void]$wsConso.Range("B2:B26").Select()

[void]$wsPro.Select()
[void]$wsPro.Range("B2:B26").Select()
